# Rework data instantiation process:
# Remove the obsolete A=60/61/62/63 parameter blocks (rows 157-176) from the
# "feasgen_thermal" sheet - the underlying data for those l1..l5 categories was
# superseded by the A=70 / 70.5 / 100 blocks that already followed them, which
# shift up into their place. Also reflect the resulting change of active
# sheet/selection that came along with this edit in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("feasgen_thermal")

# Drop the twenty now-obsolete rows (four 5-row blocks for A = 60, 61, 62, 63).
# Excel shifts everything below up automatically, so what used to be rows
# 177-191 (A = 70, 70.5, 100) becomes rows 157-171.
$ws.Rows("157:176").Delete()

# This sheet became the active / focused sheet in the saved workbook.
$ws.Activate()
$ws.Range("F153").Select()
